$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2022_9")

$ws.Cells.Item(288, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(288, 2).Value = "09:32:23 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(288, 3).Value = "developer@nex-softwares.com"
$ws.Cells.Item(288, 4).Value = "User"
$ws.Cells.Item(288, 5).Value = "/api/auth/login"
$ws.Cells.Item(288, 6).Value = "login"
$ws.Cells.Item(288, 7).Value = "succeeded"
$ws.Cells.Item(288, 8).Value = "developer@nex-softwares.com  login"

$ws.Cells.Item(289, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(289, 2).Value = "09:32:23 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(289, 4).Value = "User"
$ws.Cells.Item(289, 5).Value = "/api/user"
$ws.Cells.Item(289, 6).Value = "read"
$ws.Cells.Item(289, 7).Value = "succeeded"
$ws.Cells.Item(289, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(289, 9).Value = 1
$ws.Cells.Item(289, 10).Value = "NEX"
$ws.Cells.Item(289, 11).Value = "Admin"
$ws.Cells.Item(289, 12).Value = "all"

$ws.Cells.Item(290, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(290, 2).Value = "09:32:29 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(290, 4).Value = "User"
$ws.Cells.Item(290, 5).Value = "/api/user"
$ws.Cells.Item(290, 6).Value = "read"
$ws.Cells.Item(290, 7).Value = "succeeded"
$ws.Cells.Item(290, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(290, 9).Value = 1
$ws.Cells.Item(290, 10).Value = "NEX"
$ws.Cells.Item(290, 11).Value = "Admin"
$ws.Cells.Item(290, 12).Value = "all"

$ws.Cells.Item(291, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(291, 2).Value = "09:33:53 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(291, 4).Value = "User"
$ws.Cells.Item(291, 5).Value = "/api/user/:id"
$ws.Cells.Item(291, 6).Value = "read"
$ws.Cells.Item(291, 7).Value = "succeeded"
$ws.Cells.Item(291, 8).Value = "NEX  Admin  read user 2"
$ws.Cells.Item(291, 9).Value = 1
$ws.Cells.Item(291, 10).Value = "NEX"
$ws.Cells.Item(291, 11).Value = "Admin"
$ws.Cells.Item(291, 12).Value = "'2"
$ws.Cells.Item(291, 12).ClearFormats()

$ws.Cells.Item(292, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(292, 2).Value = "09:35:47 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(292, 4).Value = "User"
$ws.Cells.Item(292, 5).Value = "/api/user"
$ws.Cells.Item(292, 6).Value = "read"
$ws.Cells.Item(292, 7).Value = "succeeded"
$ws.Cells.Item(292, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(292, 9).Value = 1
$ws.Cells.Item(292, 10).Value = "NEX"
$ws.Cells.Item(292, 11).Value = "Admin"
$ws.Cells.Item(292, 12).Value = "all"

$ws.Cells.Item(293, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(293, 2).Value = "09:35:47 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(293, 4).Value = "User"
$ws.Cells.Item(293, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(293, 6).Value = "read"
$ws.Cells.Item(293, 7).Value = "succeeded"
$ws.Cells.Item(293, 8).Value = "NEX  Admin  read all to validate users (0) from 0 to 100"
$ws.Cells.Item(293, 9).Value = 1
$ws.Cells.Item(293, 10).Value = "NEX"
$ws.Cells.Item(293, 11).Value = "Admin"
$ws.Cells.Item(293, 12).Value = "all"

$ws.Cells.Item(294, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(294, 2).Value = "09:35:53 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(294, 4).Value = "Trip"
$ws.Cells.Item(294, 5).Value = "/api/trip"
$ws.Cells.Item(294, 6).Value = "read"
$ws.Cells.Item(294, 7).Value = "succeeded"
$ws.Cells.Item(294, 8).Value = "NEX  Admin  read all trips (1) from 0 to 100"
$ws.Cells.Item(294, 9).Value = 1
$ws.Cells.Item(294, 10).Value = "NEX"
$ws.Cells.Item(294, 11).Value = "Admin"
$ws.Cells.Item(294, 12).Value = "all"

$ws.Cells.Item(295, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(295, 2).Value = "09:37:25 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(295, 4).Value = "User"
$ws.Cells.Item(295, 5).Value = "/api/user"
$ws.Cells.Item(295, 6).Value = "read"
$ws.Cells.Item(295, 7).Value = "succeeded"
$ws.Cells.Item(295, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(295, 9).Value = 1
$ws.Cells.Item(295, 10).Value = "NEX"
$ws.Cells.Item(295, 11).Value = "Admin"
$ws.Cells.Item(295, 12).Value = "all"

$ws.Cells.Item(296, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(296, 2).Value = "09:37:26 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(296, 4).Value = "User"
$ws.Cells.Item(296, 5).Value = "/api/user"
$ws.Cells.Item(296, 6).Value = "read"
$ws.Cells.Item(296, 7).Value = "succeeded"
$ws.Cells.Item(296, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(296, 9).Value = 1
$ws.Cells.Item(296, 10).Value = "NEX"
$ws.Cells.Item(296, 11).Value = "Admin"
$ws.Cells.Item(296, 12).Value = "all"

$ws.Cells.Item(297, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(297, 2).Value = "09:37:27 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(297, 4).Value = "User"
$ws.Cells.Item(297, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(297, 6).Value = "read"
$ws.Cells.Item(297, 7).Value = "succeeded"
$ws.Cells.Item(297, 8).Value = "NEX  Admin  read all to validate users (0) from 0 to 100"
$ws.Cells.Item(297, 9).Value = 1
$ws.Cells.Item(297, 10).Value = "NEX"
$ws.Cells.Item(297, 11).Value = "Admin"
$ws.Cells.Item(297, 12).Value = "all"

$ws.Cells.Item(298, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(298, 2).Value = "09:37:28 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(298, 4).Value = "Trip"
$ws.Cells.Item(298, 5).Value = "/api/trip"
$ws.Cells.Item(298, 6).Value = "read"
$ws.Cells.Item(298, 7).Value = "succeeded"
$ws.Cells.Item(298, 8).Value = "NEX  Admin  read all trips (1) from 0 to 100"
$ws.Cells.Item(298, 9).Value = 1
$ws.Cells.Item(298, 10).Value = "NEX"
$ws.Cells.Item(298, 11).Value = "Admin"
$ws.Cells.Item(298, 12).Value = "all"

$ws.Cells.Item(299, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(299, 2).Value = "09:37:48 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(299, 4).Value = "User"
$ws.Cells.Item(299, 5).Value = "/api/user"
$ws.Cells.Item(299, 6).Value = "read"
$ws.Cells.Item(299, 7).Value = "succeeded"
$ws.Cells.Item(299, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(299, 9).Value = 1
$ws.Cells.Item(299, 10).Value = "NEX"
$ws.Cells.Item(299, 11).Value = "Admin"
$ws.Cells.Item(299, 12).Value = "all"

$ws.Cells.Item(300, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(300, 2).Value = "09:37:50 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(300, 4).Value = "User"
$ws.Cells.Item(300, 5).Value = "/api/user"
$ws.Cells.Item(300, 6).Value = "read"
$ws.Cells.Item(300, 7).Value = "succeeded"
$ws.Cells.Item(300, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(300, 9).Value = 1
$ws.Cells.Item(300, 10).Value = "NEX"
$ws.Cells.Item(300, 11).Value = "Admin"
$ws.Cells.Item(300, 12).Value = "all"

$ws.Cells.Item(301, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(301, 2).Value = "09:37:52 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(301, 4).Value = "User"
$ws.Cells.Item(301, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(301, 6).Value = "read"
$ws.Cells.Item(301, 7).Value = "succeeded"
$ws.Cells.Item(301, 8).Value = "NEX  Admin  read all to validate users (0) from 0 to 100"
$ws.Cells.Item(301, 9).Value = 1
$ws.Cells.Item(301, 10).Value = "NEX"
$ws.Cells.Item(301, 11).Value = "Admin"
$ws.Cells.Item(301, 12).Value = "all"

$ws.Cells.Item(302, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(302, 2).Value = "09:37:52 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(302, 4).Value = "Trip"
$ws.Cells.Item(302, 5).Value = "/api/trip"
$ws.Cells.Item(302, 6).Value = "read"
$ws.Cells.Item(302, 7).Value = "succeeded"
$ws.Cells.Item(302, 8).Value = "NEX  Admin  read all trips (1) from 0 to 100"
$ws.Cells.Item(302, 9).Value = 1
$ws.Cells.Item(302, 10).Value = "NEX"
$ws.Cells.Item(302, 11).Value = "Admin"
$ws.Cells.Item(302, 12).Value = "all"

$ws.Cells.Item(303, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(303, 2).Value = "09:38:10 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(303, 4).Value = "Preference"
$ws.Cells.Item(303, 5).Value = "/api/preference"
$ws.Cells.Item(303, 6).Value = "read"
$ws.Cells.Item(303, 7).Value = "succeeded"
$ws.Cells.Item(303, 8).Value = "NEX  Admin  read all preferences (undefined) from undefined to NaN"
$ws.Cells.Item(303, 9).Value = 1
$ws.Cells.Item(303, 10).Value = "NEX"
$ws.Cells.Item(303, 11).Value = "Admin"
$ws.Cells.Item(303, 12).Value = "all"

$ws.Cells.Item(304, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(304, 2).Value = "09:39:17 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(304, 4).Value = "Preference"
$ws.Cells.Item(304, 5).Value = "/api/preference"
$ws.Cells.Item(304, 6).Value = "write"
$ws.Cells.Item(304, 7).Value = "succeeded"
$ws.Cells.Item(304, 8).Value = "NEX  Admin  write preference 2"
$ws.Cells.Item(304, 9).Value = 1
$ws.Cells.Item(304, 10).Value = "NEX"
$ws.Cells.Item(304, 11).Value = "Admin"
$ws.Cells.Item(304, 12).Value = 2

$ws.Cells.Item(305, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(305, 2).Value = "09:39:18 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(305, 4).Value = "Preference"
$ws.Cells.Item(305, 5).Value = "/api/preference"
$ws.Cells.Item(305, 6).Value = "read"
$ws.Cells.Item(305, 7).Value = "succeeded"
$ws.Cells.Item(305, 8).Value = "NEX  Admin  read all preferences (undefined) from undefined to NaN"
$ws.Cells.Item(305, 9).Value = 1
$ws.Cells.Item(305, 10).Value = "NEX"
$ws.Cells.Item(305, 11).Value = "Admin"
$ws.Cells.Item(305, 12).Value = "all"

$ws.Cells.Item(306, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(306, 2).Value = "09:42:04 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(306, 4).Value = "Pricing"
$ws.Cells.Item(306, 5).Value = "/api/pricing"
$ws.Cells.Item(306, 6).Value = "read"
$ws.Cells.Item(306, 7).Value = "succeeded"
$ws.Cells.Item(306, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(306, 9).Value = 1
$ws.Cells.Item(306, 10).Value = "NEX"
$ws.Cells.Item(306, 11).Value = "Admin"
$ws.Cells.Item(306, 12).Value = "all"

$ws.Cells.Item(307, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(307, 2).Value = "09:42:18 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(307, 4).Value = "VehicleType"
$ws.Cells.Item(307, 5).Value = "/api/vehicle-type"
$ws.Cells.Item(307, 6).Value = "read"
$ws.Cells.Item(307, 7).Value = "succeeded"
$ws.Cells.Item(307, 8).Value = "NEX  Admin  read all vehicle types (undefined) from undefined to NaN"
$ws.Cells.Item(307, 9).Value = 1
$ws.Cells.Item(307, 10).Value = "NEX"
$ws.Cells.Item(307, 11).Value = "Admin"
$ws.Cells.Item(307, 12).Value = "all"

$ws.Cells.Item(308, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(308, 2).Value = "09:44:17 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(308, 4).Value = "Pricing"
$ws.Cells.Item(308, 5).Value = "/api/pricing"
$ws.Cells.Item(308, 6).Value = "read"
$ws.Cells.Item(308, 7).Value = "succeeded"
$ws.Cells.Item(308, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(308, 9).Value = 1
$ws.Cells.Item(308, 10).Value = "NEX"
$ws.Cells.Item(308, 11).Value = "Admin"
$ws.Cells.Item(308, 12).Value = "all"

$ws.Cells.Item(309, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(309, 2).Value = "09:46:45 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(309, 4).Value = "User"
$ws.Cells.Item(309, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(309, 6).Value = "read"
$ws.Cells.Item(309, 7).Value = "succeeded"
$ws.Cells.Item(309, 8).Value = "NEX  Admin  read all to validate users (0) from 0 to 100"
$ws.Cells.Item(309, 9).Value = 1
$ws.Cells.Item(309, 10).Value = "NEX"
$ws.Cells.Item(309, 11).Value = "Admin"
$ws.Cells.Item(309, 12).Value = "all"

$ws.Cells.Item(310, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(310, 2).Value = "09:48:53 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(310, 4).Value = "User"
$ws.Cells.Item(310, 5).Value = "/api/auth/upload-document"
$ws.Cells.Item(310, 6).Value = "upload"
$ws.Cells.Item(310, 7).Value = "succeeded"
$ws.Cells.Item(310, 8).Value = "John  DOE  upload his documents"
$ws.Cells.Item(310, 9).Value = 2
$ws.Cells.Item(310, 10).Value = "John"
$ws.Cells.Item(310, 11).Value = "DOE"
$ws.Cells.Item(310, 12).Value = "'2"
$ws.Cells.Item(310, 12).ClearFormats()

$ws.Cells.Item(311, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(311, 2).Value = "09:49:02 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(311, 4).Value = "User"
$ws.Cells.Item(311, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(311, 6).Value = "read"
$ws.Cells.Item(311, 7).Value = "succeeded"
$ws.Cells.Item(311, 8).Value = "NEX  Admin  read all to validate users (1) from 0 to 100"
$ws.Cells.Item(311, 9).Value = 1
$ws.Cells.Item(311, 10).Value = "NEX"
$ws.Cells.Item(311, 11).Value = "Admin"
$ws.Cells.Item(311, 12).Value = "all"

$ws.Cells.Item(312, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(312, 2).Value = "10:06:23 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(312, 4).Value = "User"
$ws.Cells.Item(312, 5).Value = "/api/user"
$ws.Cells.Item(312, 6).Value = "read"
$ws.Cells.Item(312, 7).Value = "succeeded"
$ws.Cells.Item(312, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(312, 9).Value = 1
$ws.Cells.Item(312, 10).Value = "NEX"
$ws.Cells.Item(312, 11).Value = "Admin"
$ws.Cells.Item(312, 12).Value = "all"

$ws.Cells.Item(313, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(313, 2).Value = "10:06:57 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(313, 4).Value = "User"
$ws.Cells.Item(313, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(313, 6).Value = "read"
$ws.Cells.Item(313, 7).Value = "succeeded"
$ws.Cells.Item(313, 8).Value = "NEX  Admin  read all to validate users (1) from 0 to 100"
$ws.Cells.Item(313, 9).Value = 1
$ws.Cells.Item(313, 10).Value = "NEX"
$ws.Cells.Item(313, 11).Value = "Admin"
$ws.Cells.Item(313, 12).Value = "all"

$ws.Cells.Item(314, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(314, 2).Value = "10:06:58 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(314, 4).Value = "Trip"
$ws.Cells.Item(314, 5).Value = "/api/trip"
$ws.Cells.Item(314, 6).Value = "read"
$ws.Cells.Item(314, 7).Value = "succeeded"
$ws.Cells.Item(314, 8).Value = "NEX  Admin  read all trips (1) from 0 to 100"
$ws.Cells.Item(314, 9).Value = 1
$ws.Cells.Item(314, 10).Value = "NEX"
$ws.Cells.Item(314, 11).Value = "Admin"
$ws.Cells.Item(314, 12).Value = "all"

$ws.Cells.Item(315, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(315, 2).Value = "10:07:02 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(315, 4).Value = "Preference"
$ws.Cells.Item(315, 5).Value = "/api/preference"
$ws.Cells.Item(315, 6).Value = "read"
$ws.Cells.Item(315, 7).Value = "succeeded"
$ws.Cells.Item(315, 8).Value = "NEX  Admin  read all preferences (undefined) from undefined to NaN"
$ws.Cells.Item(315, 9).Value = 1
$ws.Cells.Item(315, 10).Value = "NEX"
$ws.Cells.Item(315, 11).Value = "Admin"
$ws.Cells.Item(315, 12).Value = "all"

$ws.Cells.Item(316, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(316, 2).Value = "10:07:06 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(316, 4).Value = "VehicleType"
$ws.Cells.Item(316, 5).Value = "/api/vehicle-type"
$ws.Cells.Item(316, 6).Value = "read"
$ws.Cells.Item(316, 7).Value = "succeeded"
$ws.Cells.Item(316, 8).Value = "NEX  Admin  read all vehicle types (undefined) from undefined to NaN"
$ws.Cells.Item(316, 9).Value = 1
$ws.Cells.Item(316, 10).Value = "NEX"
$ws.Cells.Item(316, 11).Value = "Admin"
$ws.Cells.Item(316, 12).Value = "all"

$ws.Cells.Item(317, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(317, 2).Value = "10:07:08 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(317, 4).Value = "Pricing"
$ws.Cells.Item(317, 5).Value = "/api/pricing"
$ws.Cells.Item(317, 6).Value = "read"
$ws.Cells.Item(317, 7).Value = "succeeded"
$ws.Cells.Item(317, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(317, 9).Value = 1
$ws.Cells.Item(317, 10).Value = "NEX"
$ws.Cells.Item(317, 11).Value = "Admin"
$ws.Cells.Item(317, 12).Value = "all"

$ws.Cells.Item(318, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(318, 2).Value = "10:07:17 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(318, 4).Value = "User"
$ws.Cells.Item(318, 5).Value = "/api/user"
$ws.Cells.Item(318, 6).Value = "read"
$ws.Cells.Item(318, 7).Value = "succeeded"
$ws.Cells.Item(318, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(318, 9).Value = 1
$ws.Cells.Item(318, 10).Value = "NEX"
$ws.Cells.Item(318, 11).Value = "Admin"
$ws.Cells.Item(318, 12).Value = "all"

$ws.Cells.Item(319, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(319, 2).Value = "10:11:03 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(319, 4).Value = "User"
$ws.Cells.Item(319, 5).Value = "/api/auth/logout"
$ws.Cells.Item(319, 6).Value = "logout"
$ws.Cells.Item(319, 7).Value = "succeeded"
$ws.Cells.Item(319, 8).Value = "NEX  Admin  logout "
$ws.Cells.Item(319, 9).Value = 1
$ws.Cells.Item(319, 10).Value = "NEX"
$ws.Cells.Item(319, 11).Value = "Admin"
$ws.Cells.Item(319, 12).Value = "'1"
$ws.Cells.Item(319, 12).ClearFormats()

$ws.Cells.Item(320, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(320, 2).Value = "13:24:35 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(320, 3).Value = "developer@nex-softwares.com"
$ws.Cells.Item(320, 4).Value = "User"
$ws.Cells.Item(320, 5).Value = "/api/auth/login"
$ws.Cells.Item(320, 6).Value = "login"
$ws.Cells.Item(320, 7).Value = "succeeded"
$ws.Cells.Item(320, 8).Value = "developer@nex-softwares.com  login"

$ws.Cells.Item(321, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(321, 2).Value = "13:24:36 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(321, 4).Value = "User"
$ws.Cells.Item(321, 5).Value = "/api/user"
$ws.Cells.Item(321, 6).Value = "read"
$ws.Cells.Item(321, 7).Value = "succeeded"
$ws.Cells.Item(321, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(321, 9).Value = 1
$ws.Cells.Item(321, 10).Value = "NEX"
$ws.Cells.Item(321, 11).Value = "Admin"
$ws.Cells.Item(321, 12).Value = "all"

$ws.Cells.Item(322, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(322, 2).Value = "13:24:39 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(322, 4).Value = "User"
$ws.Cells.Item(322, 5).Value = "/api/user"
$ws.Cells.Item(322, 6).Value = "read"
$ws.Cells.Item(322, 7).Value = "succeeded"
$ws.Cells.Item(322, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(322, 9).Value = 1
$ws.Cells.Item(322, 10).Value = "NEX"
$ws.Cells.Item(322, 11).Value = "Admin"
$ws.Cells.Item(322, 12).Value = "all"

$ws.Cells.Item(323, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(323, 2).Value = "13:24:52 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(323, 4).Value = "User"
$ws.Cells.Item(323, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(323, 6).Value = "read"
$ws.Cells.Item(323, 7).Value = "succeeded"
$ws.Cells.Item(323, 8).Value = "NEX  Admin  read all to validate users (1) from 0 to 100"
$ws.Cells.Item(323, 9).Value = 1
$ws.Cells.Item(323, 10).Value = "NEX"
$ws.Cells.Item(323, 11).Value = "Admin"
$ws.Cells.Item(323, 12).Value = "all"

$ws.Cells.Item(324, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(324, 2).Value = "13:25:25 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(324, 4).Value = "User"
$ws.Cells.Item(324, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(324, 6).Value = "read"
$ws.Cells.Item(324, 7).Value = "succeeded"
$ws.Cells.Item(324, 8).Value = "NEX  Admin  read all to validate users (1) from 0 to 100"
$ws.Cells.Item(324, 9).Value = 1
$ws.Cells.Item(324, 10).Value = "NEX"
$ws.Cells.Item(324, 11).Value = "Admin"
$ws.Cells.Item(324, 12).Value = "all"

$ws.Cells.Item(325, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(325, 2).Value = "13:25:56 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(325, 4).Value = "Trip"
$ws.Cells.Item(325, 5).Value = "/api/trip"
$ws.Cells.Item(325, 6).Value = "read"
$ws.Cells.Item(325, 7).Value = "succeeded"
$ws.Cells.Item(325, 8).Value = "NEX  Admin  read all trips (1) from 0 to 100"
$ws.Cells.Item(325, 9).Value = 1
$ws.Cells.Item(325, 10).Value = "NEX"
$ws.Cells.Item(325, 11).Value = "Admin"
$ws.Cells.Item(325, 12).Value = "all"

$ws.Cells.Item(326, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(326, 2).Value = "13:26:19 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(326, 4).Value = "User"
$ws.Cells.Item(326, 5).Value = "/api/user"
$ws.Cells.Item(326, 6).Value = "read"
$ws.Cells.Item(326, 7).Value = "succeeded"
$ws.Cells.Item(326, 8).Value = "NEX  Admin  read all users (3) from 0 to 100"
$ws.Cells.Item(326, 9).Value = 1
$ws.Cells.Item(326, 10).Value = "NEX"
$ws.Cells.Item(326, 11).Value = "Admin"
$ws.Cells.Item(326, 12).Value = "all"

$ws.Cells.Item(327, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(327, 2).Value = "13:26:21 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(327, 4).Value = "User"
$ws.Cells.Item(327, 5).Value = "/api/user/:id"
$ws.Cells.Item(327, 6).Value = "read"
$ws.Cells.Item(327, 7).Value = "succeeded"
$ws.Cells.Item(327, 8).Value = "NEX  Admin  read user 2"
$ws.Cells.Item(327, 9).Value = 1
$ws.Cells.Item(327, 10).Value = "NEX"
$ws.Cells.Item(327, 11).Value = "Admin"
$ws.Cells.Item(327, 12).Value = "'2"
$ws.Cells.Item(327, 12).ClearFormats()

$ws.Cells.Item(328, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(328, 2).Value = "13:27:08 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(328, 4).Value = "Trip"
$ws.Cells.Item(328, 5).Value = "/api/trip"
$ws.Cells.Item(328, 6).Value = "read"
$ws.Cells.Item(328, 7).Value = "succeeded"
$ws.Cells.Item(328, 8).Value = "NEX  Admin  read all trips (1) from 0 to 100"
$ws.Cells.Item(328, 9).Value = 1
$ws.Cells.Item(328, 10).Value = "NEX"
$ws.Cells.Item(328, 11).Value = "Admin"
$ws.Cells.Item(328, 12).Value = "all"

$ws.Cells.Item(329, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(329, 2).Value = "13:27:24 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(329, 4).Value = "Preference"
$ws.Cells.Item(329, 5).Value = "/api/preference"
$ws.Cells.Item(329, 6).Value = "read"
$ws.Cells.Item(329, 7).Value = "succeeded"
$ws.Cells.Item(329, 8).Value = "NEX  Admin  read all preferences (undefined) from undefined to NaN"
$ws.Cells.Item(329, 9).Value = 1
$ws.Cells.Item(329, 10).Value = "NEX"
$ws.Cells.Item(329, 11).Value = "Admin"
$ws.Cells.Item(329, 12).Value = "all"

$ws.Cells.Item(330, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(330, 2).Value = "13:27:40 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(330, 4).Value = "VehicleType"
$ws.Cells.Item(330, 5).Value = "/api/vehicle-type"
$ws.Cells.Item(330, 6).Value = "read"
$ws.Cells.Item(330, 7).Value = "succeeded"
$ws.Cells.Item(330, 8).Value = "NEX  Admin  read all vehicle types (undefined) from undefined to NaN"
$ws.Cells.Item(330, 9).Value = 1
$ws.Cells.Item(330, 10).Value = "NEX"
$ws.Cells.Item(330, 11).Value = "Admin"
$ws.Cells.Item(330, 12).Value = "all"

$ws.Cells.Item(331, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(331, 2).Value = "13:28:00 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(331, 4).Value = "Pricing"
$ws.Cells.Item(331, 5).Value = "/api/pricing"
$ws.Cells.Item(331, 6).Value = "read"
$ws.Cells.Item(331, 7).Value = "succeeded"
$ws.Cells.Item(331, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(331, 9).Value = 1
$ws.Cells.Item(331, 10).Value = "NEX"
$ws.Cells.Item(331, 11).Value = "Admin"
$ws.Cells.Item(331, 12).Value = "all"

$ws.Cells.Item(332, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(332, 2).Value = "13:28:26 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(332, 4).Value = "Pricing"
$ws.Cells.Item(332, 5).Value = "/api/pricing"
$ws.Cells.Item(332, 6).Value = "write"
$ws.Cells.Item(332, 7).Value = "succeeded"
$ws.Cells.Item(332, 8).Value = "NEX  Admin  write pricing 8"
$ws.Cells.Item(332, 9).Value = 1
$ws.Cells.Item(332, 10).Value = "NEX"
$ws.Cells.Item(332, 11).Value = "Admin"
$ws.Cells.Item(332, 12).Value = 8

$ws.Cells.Item(333, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(333, 2).Value = "13:28:26 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(333, 4).Value = "Pricing"
$ws.Cells.Item(333, 5).Value = "/api/pricing"
$ws.Cells.Item(333, 6).Value = "read"
$ws.Cells.Item(333, 7).Value = "succeeded"
$ws.Cells.Item(333, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(333, 9).Value = 1
$ws.Cells.Item(333, 10).Value = "NEX"
$ws.Cells.Item(333, 11).Value = "Admin"
$ws.Cells.Item(333, 12).Value = "all"

$ws.Cells.Item(334, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(334, 2).Value = "13:28:44 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(334, 4).Value = "Pricing"
$ws.Cells.Item(334, 5).Value = "/api/pricing"
$ws.Cells.Item(334, 6).Value = "edit"
$ws.Cells.Item(334, 7).Value = "succeeded"
$ws.Cells.Item(334, 8).Value = "NEX  Admin  edit pricing 8"
$ws.Cells.Item(334, 9).Value = 1
$ws.Cells.Item(334, 10).Value = "NEX"
$ws.Cells.Item(334, 11).Value = "Admin"
$ws.Cells.Item(334, 12).Value = 8

$ws.Cells.Item(335, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(335, 2).Value = "13:28:44 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(335, 4).Value = "Pricing"
$ws.Cells.Item(335, 5).Value = "/api/pricing"
$ws.Cells.Item(335, 6).Value = "read"
$ws.Cells.Item(335, 7).Value = "succeeded"
$ws.Cells.Item(335, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(335, 9).Value = 1
$ws.Cells.Item(335, 10).Value = "NEX"
$ws.Cells.Item(335, 11).Value = "Admin"
$ws.Cells.Item(335, 12).Value = "all"

$ws.Cells.Item(336, 1).Value = "Sun Sep 04 2022"
$ws.Cells.Item(336, 2).Value = "13:28:51 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(336, 4).Value = "Pricing"
$ws.Cells.Item(336, 5).Value = "/api/pricing"
$ws.Cells.Item(336, 6).Value = "read"
$ws.Cells.Item(336, 7).Value = "succeeded"
$ws.Cells.Item(336, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(336, 9).Value = 1
$ws.Cells.Item(336, 10).Value = "NEX"
$ws.Cells.Item(336, 11).Value = "Admin"
$ws.Cells.Item(336, 12).Value = "all"
